# Update cryptos list snapshot — GitHub Actions scheduled refresh.
# Column D ("Price") values are textual (may contain thousands separators
# like "28.030.73" or leading zeros like "0.000008880"), so each is entered
# with a leading apostrophe to force a Text quote-prefix and keep Excel
# from re-interpreting them as numbers / dates / losing precision.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.030.73"
$ws.Range("E2").Value = "  +2.09%  "
$ws.Range("E3").Value = "  +2.70%  "
$ws.Range("E4").Value = "  -0.51%  "
$ws.Range("D5").Value = "'315.65"
$ws.Range("E5").Value = "  +1.41%  "
$ws.Range("E6").Value = "  -0.50%  "
$ws.Range("D7").Value = "'0.4813"
$ws.Range("E7").Value = "  +0.79%  "
$ws.Range("D8").Value = "'0.3817"
$ws.Range("E8").Value = "  +0.38%  "
$ws.Range("D9").Value = "'0.07362"
$ws.Range("E9").Value = "  +0.44%  "
$ws.Range("D10").Value = "'0.9336"
$ws.Range("E10").Value = "  -0.15%  "
$ws.Range("D11").Value = "'20.80"
$ws.Range("E11").Value = "  +0.29%  "
$ws.Range("D12").Value = "'0.07806"
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").Value = "'1.875.14"
$ws.Range("E13").Value = "  +0.33%  "
$ws.Range("D14").Value = "'5.507"
$ws.Range("E14").Value = "  +1.22%  "
$ws.Range("D15").Value = "'6.634"
$ws.Range("E15").Value = "  +1.24%  "
$ws.Range("D16").Value = "'91.89"
$ws.Range("E16").Value = "  +1.65%  "
$ws.Range("E17").Value = "  -0.51%  "
$ws.Range("D18").Value = "'0.000008880"
$ws.Range("E18").Value = "  +0.95%  "
$ws.Range("D19").Value = "'1.005"
$ws.Range("E19").Value = "  -0.39%  "
$ws.Range("D20").Value = "'28.059.29"
$ws.Range("E20").Value = "  +1.88%  "
$ws.Range("D21").Value = "'14.76"
$ws.Range("E21").Value = "  +0.68%  "
$ws.Range("D22").Value = "'5.176"
$ws.Range("E22").Value = "  +1.05%  "
$ws.Range("D23").Value = "'2.129.67"
$ws.Range("E23").Value = "  +1.35%  "
$ws.Range("E24").Value = "  +1.95%  "
$ws.Range("D25").Value = "'156.62"
$ws.Range("E25").Value = "  +0.93%  "
$ws.Range("D26").Value = "'1.913"
$ws.Range("E26").Value = "  -1.20%  "
$ws.Range("D27").Value = "'18.51"
$ws.Range("E27").Value = "  +0.23%  "
$ws.Range("D28").Value = "'2.131"
$ws.Range("E28").Value = "  +5.51%  "
$ws.Range("D29").Value = "'116.82"
$ws.Range("E29").Value = "  +1.30%  "
$ws.Range("D30").Value = "'4.968"
$ws.Range("E30").Value = "  +0.69%  "
$ws.Range("D31").Value = "'0.08951"
$ws.Range("E31").Value = "  +0.72%  "
$ws.Range("D32").Value = "'3.291"
$ws.Range("E32").Value = "  -0.99%  "
$ws.Range("D33").Value = "'1.255"
$ws.Range("E33").Value = "  +3.38%  "
$ws.Range("D34").Value = "'0.7755"
$ws.Range("E34").Value = "  +2.36%  "
$ws.Range("D35").Value = "'4.667"
$ws.Range("E35").Value = "  +1.52%  "
$ws.Range("D36").Value = "'2.609"
$ws.Range("E36").Value = "  -4.38%  "
$ws.Range("D37").Value = "'0.02049"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("E38").Value = "  -1.27%  "
$ws.Range("D39").Value = "'0.5519"
$ws.Range("E39").Value = "  -1.25%  "
$ws.Range("D40").Value = "'0.05305"
$ws.Range("E40").Value = "  +0.57%  "
$ws.Range("D41").Value = "'2.995"
$ws.Range("E41").Value = "  +0.14%  "
$ws.Range("D42").Value = "'7.026"
$ws.Range("E42").Value = "  -0.51%  "

# Rows 43/44 swap rank order (Aptos now ahead of Algorand) along with
# their own updated price / volume figures.
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "'8.503"
$ws.Range("E43").Value = "  -1.60%  "
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").Value = "'0.1526"
$ws.Range("E44").Value = "  +0.12%  "

$ws.Range("D45").Value = "'10.68"
$ws.Range("E45").Value = "  +0.10%  "
$ws.Range("D46").Value = "'108.78"
$ws.Range("E46").Value = "  +5.68%  "
$ws.Range("D47").Value = "'0.4832"
$ws.Range("E47").Value = "  -1.66%  "
$ws.Range("D48").Value = "'1.006"
$ws.Range("E48").Value = "  -0.51%  "
$ws.Range("D49").Value = "'1.650"
$ws.Range("E49").Value = "  -0.38%  "
$ws.Range("D50").Value = "'68.01"
$ws.Range("E50").Value = "  +0.88%  "
$ws.Range("D51").Value = "'0.06079"
